$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '61.815.72'
$ws.Range('E2').Value = '  -2.26%  '
$ws.Range('D3').Value = '2.576.54'
$ws.Range('E3').Value = '  -4.31%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '548.15'
$ws.Range('E5').Value = '  -1.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '154.48'
$ws.Range('E6').Value = '  -1.40%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('E8').Value = '  +2.54%  '
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('E10').Value = '  -1.61%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.364'
$ws.Range('E12').Value = '  -1.05%  '
$ws.Range('D13').Value = '3.029.42'
$ws.Range('E13').Value = '  -4.57%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '25.41'
$ws.Range('E14').Value = '  -3.41%  '
$ws.Range('D15').Value = '61.704.10'
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('E16').Value = '  -0.57%  '
$ws.Range('D17').Value = '2.578.09'
$ws.Range('E17').Value = '  -4.71%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '11.51'
$ws.Range('E18').Value = '  -4.34%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '4.54'
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '337.19'
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('E21').Value = '  -4.61%  '
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.492'
$ws.Range('E23').Value = '  -2.45%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '63.31'
$ws.Range('E24').Value = '  -0.64%  '
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '8.07'
$ws.Range('E27').Value = '  +0.29%  '
$ws.Range('E28').Value = '  +6.18%  '
$ws.Range('D29').Value = '0.0₃0838'
$ws.Range('E29').Value = '  -2.49%  '
$ws.Range('E30').Value = '  -1.75%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.29'
$ws.Range('E31').Value = '  -2.67%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '160.73'
$ws.Range('E32').Value = '  -2.68%  '
$ws.Range('E33').Value = '  +0.02%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '4.74'
$ws.Range('E34').Value = '  -0.35%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '19.15'
$ws.Range('E35').Value = '  -1.92%  '
$ws.Range('E36').Value = '  -1.91%  '
$ws.Range('E37').Value = '  +0.59%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '332.85'
$ws.Range('E38').Value = '  -1.96%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.918'
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '5.89'
$ws.Range('E40').Value = '  -2.63%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '3.94'
$ws.Range('E41').Value = '  -0.04%  '
$ws.Range('E42').Value = '  -1.48%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '20.69'
$ws.Range('E43').Value = '  -0.74%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.997'
$ws.Range('E44').Value = '  -0.17%  '
$ws.Range('D45').Value = '2.127.84'
$ws.Range('E45').Value = '  +0.95%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.605'
$ws.Range('E46').Value = '  -2.48%  '
$ws.Range('E47').Value = '  -1.12%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '19.56'
$ws.Range('E48').Value = '  -3.55%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0545'
$ws.Range('E49').Value = '  -2.95%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0965'
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('E51').Value = '  -1.29%  '
